$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-95 down to 27-96.
$ws.Rows.Item(26).Insert()

# The inherited formatting/merges for the "Plantation forests" secondary-category
# block (previously C22:C25 / D22:D25) must expand to include the newly
# inserted row so it keeps reading as one Secondary class.
$ws.Range("C22:C25").UnMerge()
$ws.Range("C22:C26").Merge()
$ws.Range("D22:D25").UnMerge()
$ws.Range("D22:D26").Merge()

# Populate the new tertiary-class row: "2.1.5 Permanent carbon forest"
# (Primary/Secondary code & class cells stay blank - they are covered by the
# merged cells above/around them).
$ws.Range("E26").Value = "2.1.5"
$ws.Range("F26").Value = "Permanent carbon forest"
